$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1157.4286
$ws.Range("J19").Value = 1145.3
$ws.Range("L19").Value = 1145.3
$ws.Range("N19").Value = -1495.3
$ws.Range("H32").Value = 3871.2
$ws.Range("I32").Value = 3227.2856
$ws.Range("J32").Value = 4434.625
$ws.Range("K32").Value = 3227.2856
$ws.Range("L32").Value = 4434.625
$ws.Range("M32").Value = -2901.2856
$ws.Range("N32").Value = -5086.625
$ws.Range("H40").Value = 6661.4165
$ws.Range("I40").Value = 6296.3335
$ws.Range("J40").Value = 6783.1113
$ws.Range("K40").Value = 6296.3335
$ws.Range("L40").Value = 6783.1113
$ws.Range("M40").Value = -6121.3335
$ws.Range("N40").Value = -7133.1113
$ws.Range("H41").Value = 1011
$ws.Range("I41").Value = 1122.5
$ws.Range("J41").Value = 899.5
$ws.Range("K41").Value = 1122.5
$ws.Range("L41").Value = 899.5
$ws.Range("M41").Value = -682.5
$ws.Range("N41").Value = -1779.5
$ws.Range("H55").Value = 103.933334
$ws.Range("I55").Value = 41.666668
$ws.Range("K55").Value = 41.666668
$ws.Range("M55").Value = 172.333332
$ws.Range("H58").Value = 197.16667
$ws.Range("I58").Value = 197.16667
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 591.50001
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -441.50001
$ws.Range("N58").ClearContents()
$ws.Range("H61").Value = 1149.5
$ws.Range("I61").Value = 1149.5
$ws.Range("K61").Value = 3448.5
$ws.Range("M61").Value = -3276.5
$ws.Range("H76").Value = 4999.5
$ws.Range("I76").Value = 4999
$ws.Range("K76").Value = 4999
$ws.Range("M76").Value = -4684
$ws.Range("H79").Value = 4999.5
$ws.Range("I79").Value = 4999
$ws.Range("K79").Value = 4999
$ws.Range("M79").Value = -3907
$ws.Range("H87").Value = 74999.5
$ws.Range("J87").Value = 74999.5
$ws.Range("L87").Value = 74999.5
$ws.Range("N87").Value = -77495.5
$ws.Range("H90").Value = 74999.5
$ws.Range("J90").Value = 74999.5
$ws.Range("L90").Value = 224998.5
$ws.Range("N90").Value = -237478.5
$ws.Range("H98").Value = 1103.3334
$ws.Range("I98").Value = 1163.7
$ws.Range("J98").Value = 801.5
$ws.Range("K98").Value = 1163.7
$ws.Range("L98").Value = 801.5
$ws.Range("M98").Value = 334.3
$ws.Range("N98").Value = -3797.5
$ws.Range("H122").Value = 1103.3334
$ws.Range("I122").Value = 1163.7
$ws.Range("J122").Value = 801.5
$ws.Range("K122").Value = 3491.1
$ws.Range("L122").Value = 2404.5
$ws.Range("M122").Value = -1041.1
$ws.Range("N122").Value = -7304.5
$ws.Range("H131").Value = 1985.375
$ws.Range("I131").Value = 2054.7144
$ws.Range("J131").Value = 1500
$ws.Range("K131").Value = 6164.1432
$ws.Range("L131").Value = 4500
$ws.Range("M131").Value = -1124.1432
$ws.Range("N131").Value = -14580
$ws.Range("H141").Value = 6295
$ws.Range("I141").Value = 3825
$ws.Range("K141").Value = 11475
$ws.Range("M141").Value = -6295

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 674
$ws.Range("I2").Value = 674
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 674
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -561
$ws.Range("N2").ClearContents()
$ws.Range("H116").Value = 674
$ws.Range("I116").Value = 674
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 674
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 1620
$ws.Range("N116").ClearContents()
$ws.Range("H132").Value = 3071.5676
$ws.Range("I132").Value = 3038.6897
$ws.Range("K132").Value = 9116.069100000001
$ws.Range("M132").Value = -6586.069100000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 674
$ws.Range("I3").Value = 674
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 674
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -560
$ws.Range("N3").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9691.896000000001
$ws.Range("I31").Value = 3571.8667
$ws.Range("J31").Value = 16249.071
$ws.Range("K31").Value = 3571.8667
$ws.Range("L31").Value = 16249.071
$ws.Range("M31").Value = -3276.8667
$ws.Range("N31").Value = -16839.071
$ws.Range("H34").Value = 9691.896000000001
$ws.Range("I34").Value = 3571.8667
$ws.Range("J34").Value = 16249.071
$ws.Range("K34").Value = 3571.8667
$ws.Range("L34").Value = 16249.071
$ws.Range("M34").Value = -3369.8667
$ws.Range("N34").Value = -16653.071
$ws.Range("H41").Value = 40999.8
$ws.Range("I41").Value = 30000
$ws.Range("J41").Value = 43749.75
$ws.Range("K41").Value = 30000
$ws.Range("L41").Value = 43749.75
$ws.Range("M41").Value = -29572
$ws.Range("N41").Value = -44605.75
$ws.Range("H47").Value = 26666
$ws.Range("J47").Value = 26666
$ws.Range("L47").Value = 26666
$ws.Range("N47").Value = -27798
$ws.Range("H60").Value = 35376.23
$ws.Range("I60").Value = 18749
$ws.Range("J60").Value = 38399.363
$ws.Range("K60").Value = 18749
$ws.Range("L60").Value = 38399.363
$ws.Range("M60").Value = -18238
$ws.Range("N60").Value = -39421.363
$ws.Range("H139").Value = 119999
$ws.Range("J139").Value = 119999
$ws.Range("L139").Value = 119999
$ws.Range("N139").Value = -130279
$ws.Range("H141").Value = 212355.06
$ws.Range("J141").Value = 256799.64
$ws.Range("L141").Value = 256799.64
$ws.Range("N141").Value = -267159.64

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 180.4
$ws.Range("J23").Value = 226
$ws.Range("L23").Value = 678
$ws.Range("N23").Value = -1148
$ws.Range("H61").Value = 193.63637
$ws.Range("I61").Value = 168
$ws.Range("K61").Value = 504
$ws.Range("M61").Value = -289
$ws.Range("H80").Value = 4999
$ws.Range("J80").Value = 4999
$ws.Range("L80").Value = 14997
$ws.Range("N80").Value = -16869
$ws.Range("H83").Value = 4999
$ws.Range("J83").Value = 4999
$ws.Range("L83").Value = 44991
$ws.Range("N83").Value = -54351
$ws.Range("H97").Value = 1139.1111
$ws.Range("I97").Value = 52
$ws.Range("J97").Value = 1275
$ws.Range("K97").Value = 156
$ws.Range("L97").Value = 3825
$ws.Range("M97").Value = 340
$ws.Range("N97").Value = -4817

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 5669.4443
$ws.Range("I122").Value = 1995.5454
$ws.Range("J122").Value = 11442.714
$ws.Range("K122").Value = 5986.6362
$ws.Range("L122").Value = 34328.142
$ws.Range("M122").Value = -3536.6362
$ws.Range("N122").Value = -39228.142
$ws.Range("H132").Value = 4761.1377
$ws.Range("I132").Value = 3049.9048
$ws.Range("J132").Value = 9253.125
$ws.Range("K132").Value = 9149.714399999999
$ws.Range("L132").Value = 27759.375
$ws.Range("M132").Value = -6619.714399999999
$ws.Range("N132").Value = -32819.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 7237.0625
$ws.Range("I16").Value = 1097.9286
$ws.Range("J16").Value = 50211
$ws.Range("K16").Value = 1097.9286
$ws.Range("L16").Value = 50211
$ws.Range("M16").Value = -927.9286
$ws.Range("N16").Value = -50551
$ws.Range("H40").Value = 7595.28
$ws.Range("I40").Value = 6951.5264
$ws.Range("K40").Value = 6951.5264
$ws.Range("M40").Value = -6815.5264
$ws.Range("H46").Value = 2424.4119
$ws.Range("I46").Value = 683.4
$ws.Range("K46").Value = 683.4
$ws.Range("M46").Value = -495.4
$ws.Range("H136").Value = 9088.074000000001
$ws.Range("I136").Value = 5345.4614
$ws.Range("J136").Value = 9989.074000000001
$ws.Range("K136").Value = 16036.3842
$ws.Range("L136").Value = 29967.222
$ws.Range("M136").Value = -13486.3842
$ws.Range("N136").Value = -35067.222

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4148
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 4148
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()
$ws.Range("H126").Value = 3293.2307
$ws.Range("I126").Value = 3293.2307
$ws.Range("K126").Value = 9879.6921
$ws.Range("M126").Value = -7409.6921
$ws.Range("H136").Value = 7033.706
$ws.Range("I136").Value = 5659.846
$ws.Range("K136").Value = 16979.538
$ws.Range("M136").Value = -14429.538
